# Regenerate the "K" (strikeouts) column (G) for the save_data sheet.
# The sheet's header row labels column G as "K" (previously it held a
# "Strike#" style raw count). This recalculates/rewrites the per-game
# K values (s_vals) that belong in column G for rows 3-56.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> recalculated K value
$kValues = [ordered]@{
    3  = 0
    4  = 2
    5  = 1
    6  = 0
    7  = 1
    8  = 2
    9  = 0
    11 = 0
    12 = 1
    13 = 2
    14 = 0
    15 = 0
    16 = 1
    17 = 1
    18 = 2
    19 = 2
    20 = 0
    21 = 1
    22 = 1
    23 = 1
    24 = 0
    25 = 0
    26 = 2
    27 = 1
    28 = 1
    29 = 0
    30 = 1
    31 = 0
    32 = 0
    33 = 1
    34 = 1
    35 = 1
    36 = 2
    37 = 0
    38 = 0
    39 = 0
    40 = 1
    41 = 0
    42 = 0
    43 = 2
    45 = 1
    46 = 2
    47 = 1
    48 = 0
    49 = 0
    50 = 2
    51 = 0
    52 = 2
    53 = 3
    54 = 2
    55 = 0
    56 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
